$wb = $excel.ActiveWorkbook

function Set-DataRow {
    param(
        $ws,
        [int]$row,
        [object[]]$values
    )
    # values is an array the same length as the header row (columns A..).
    for ($c = 1; $c -le $values.Length; $c++) {
        $v = $values[$c - 1]
        if ($null -ne $v) {
            $cell = $ws.Cells.Item($row, $c)
            $cell.Value = $v
            if ($c -eq 1) {
                # column A (index) cells use the bold/bordered header-like style
                $cell.Font.Bold = $true
                $cell.Borders.LineStyle = 1
                $cell.HorizontalAlignment = -4108
                $cell.VerticalAlignment = -4160
            }
        }
    }
}

# ---------------------------------------------------------------------------
# Sheet "土地" (land) - append new data row 2
# ---------------------------------------------------------------------------
$wsLand = $wb.Worksheets.Item("土地")
Set-DataRow $wsLand 2 @(
    14,
    "臺北市大安區龍泉段一小段02930000地號",
    365,
    "100000分之16216",
    "高金素梅",
    "92年12月25日",
    "033貝賣",
    "25000000(土地建物與車位合併價）",
    "land",
    "normal",
    "2012-04-30",
    "高金素梅",
    926,
    "tmp92521",
    14,
    0.16216,
    59.1884
)

# ---------------------------------------------------------------------------
# Sheet "建物" (building) - rebuild header row to the full 17-column schema
# and append new data row 2
# ---------------------------------------------------------------------------
$wsBuilding = $wb.Worksheets.Item("建物")
$headers = @(
    "name","area","share_portion","owner","register_date","register_reason",
    "acquire_value","property_category","category","date","legislator_name",
    "legislator_id","source_file","index","portion","total"
)
for ($c = 2; $c -le 17; $c++) {
    $wsBuilding.Cells.Item(1, $c).Value = $headers[$c - 2]
}

Set-DataRow $wsBuilding 2 @(
    19,
    "臺北市大安區龍泉段一小段05819000建號",
    148.31,
    "全部",
    "高金素梅",
    "92年12月25日",
    "貝賣",
    "25000000(土地建物與車位合併價）",
    "land",
    "normal",
    "2012-04-30",
    "高金素梅",
    926,
    "tmp92521",
    19,
    1,
    148.31
)

# ---------------------------------------------------------------------------
# Sheet "存款" (deposit) - append new data row 2 (duplicate of row 1, plus index)
# ---------------------------------------------------------------------------
$wsDeposit = $wb.Worksheets.Item("存款")
Set-DataRow $wsDeposit 2 @(
    45,
    "臺灣銀行群賢分行",
    "活期存款",
    "新臺幣",
    "高金素梅",
    366738
)

# ---------------------------------------------------------------------------
# Sheet "債務" (debt) - insert a new row 2 (existing row 2 shifts down to row 3)
# ---------------------------------------------------------------------------
$wsDebt = $wb.Worksheets.Item("債務")
$wsDebt.Rows.Item(2).Insert()

Set-DataRow $wsDebt 2 @(
    85,
    "現金",
    "高金素梅",
    "陳麗卿新北市泰山區明志路",
    "6000000",
    "96年02月06日",
    "借款"
)
